# Applies the "[AFG] added final excel sheets for Afghanistan" edit:
#  - Adds two new worksheets: "ODI Batting Extra" and "ODI Bowling Extra"
#  - Populates them with per-match supplementary stats
#  - Removes now-redundant empty placeholder cells from "ODI Batting" column B

$wb = $excel.ActiveWorkbook

# --- Clean up stray empty placeholder cells on the existing "ODI Batting" sheet ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$emptyBattingCells = @("B3","B7","B13","B14","B15","B16","B18","B20")
foreach ($cellRef in $emptyBattingCells) {
    [void]$wsBatting.Range($cellRef).ClearContents()
}

# --- Add "ODI Batting Extra" sheet after "ODI Bowling" ---
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBattingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBowling)
$wsBattingExtra.Name = "ODI Batting Extra"

# --- Add "ODI Bowling Extra" sheet after "ODI Batting Extra" ---
$wsBowlingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBattingExtra)
$wsBowlingExtra.Name = "ODI Bowling Extra"

# --- Page margins matching the rest of the workbook ---
foreach ($sheet in @($wsBattingExtra, $wsBowlingExtra)) {
    $sheet.PageSetup.LeftMargin = 0.75 * 72
    $sheet.PageSetup.RightMargin = 0.75 * 72
    $sheet.PageSetup.TopMargin = 1 * 72
    $sheet.PageSetup.BottomMargin = 1 * 72
    $sheet.PageSetup.HeaderMargin = 0.5 * 72
    $sheet.PageSetup.FooterMargin = 0.5 * 72
}

# --- "ODI Batting Extra" header row ---
$wsBattingExtra.Range("A1").Value = "'MATCH_CODE"
$wsBattingExtra.Range("B1").Value = "'BATTING_POSITION"
$wsBattingExtra.Range("C1").Value = "'NUM_4"
$wsBattingExtra.Range("D1").Value = "'NUM_6"
$wsBattingExtra.Range("E1").Value = "'PERCENT_RUNS_OF_TOTAL"
$wsBattingExtra.Range("F1").Value = "'MAN_OF_MATCH"

# --- "ODI Batting Extra" data rows ---
$wsBattingExtra.Range("A2").Value = "'3649"
$wsBattingExtra.Range("F2").Value = "'NO"
$wsBattingExtra.Range("A3").Value = "'3650"
$wsBattingExtra.Range("B3").Value = 10
$wsBattingExtra.Range("F3").Value = "'NO"
$wsBattingExtra.Range("A4").Value = "'3651"
$wsBattingExtra.Range("F4").Value = "'NO"
$wsBattingExtra.Range("A5").Value = "'3652"
$wsBattingExtra.Range("B5").Value = 10
$wsBattingExtra.Range("C5").Value = "'0"
$wsBattingExtra.Range("D5").Value = "'0"
$wsBattingExtra.Range("E5").Value = "'1.16%"
$wsBattingExtra.Range("F5").Value = "'NO"
$wsBattingExtra.Range("A6").Value = "'3707"
$wsBattingExtra.Range("F6").Value = "'NO"
$wsBattingExtra.Range("A7").Value = "'4134"
$wsBattingExtra.Range("B7").Value = 6
$wsBattingExtra.Range("F7").Value = "'NO"
$wsBattingExtra.Range("A8").Value = "'4136"
$wsBattingExtra.Range("B8").Value = 8
$wsBattingExtra.Range("C8").Value = "'1"
$wsBattingExtra.Range("D8").Value = "'1"
$wsBattingExtra.Range("E8").Value = "'8.71%"
$wsBattingExtra.Range("F8").Value = "'YES"
$wsBattingExtra.Range("A9").Value = "'4140"
$wsBattingExtra.Range("B9").Value = 7
$wsBattingExtra.Range("C9").Value = "'0"
$wsBattingExtra.Range("D9").Value = "'0"
$wsBattingExtra.Range("E9").Value = "'2.35%"
$wsBattingExtra.Range("F9").Value = "'NO"
$wsBattingExtra.Range("A10").Value = "'4145"
$wsBattingExtra.Range("B10").Value = 7
$wsBattingExtra.Range("C10").Value = "'0"
$wsBattingExtra.Range("D10").Value = "'0"
$wsBattingExtra.Range("F10").Value = "'NO"
$wsBattingExtra.Range("A11").Value = "'4147"
$wsBattingExtra.Range("F11").Value = "'NO"
$wsBattingExtra.Range("A12").Value = "'4154"
$wsBattingExtra.Range("F12").Value = "'NO"
$wsBattingExtra.Range("A13").Value = "'4160"
$wsBattingExtra.Range("B13").Value = 9
$wsBattingExtra.Range("F13").Value = "'NO"
$wsBattingExtra.Range("A14").Value = "'4163"
$wsBattingExtra.Range("B14").Value = 9
$wsBattingExtra.Range("F14").Value = "'NO"
$wsBattingExtra.Range("A15").Value = "'4164"
$wsBattingExtra.Range("B15").Value = 9
$wsBattingExtra.Range("F15").Value = "'NO"
$wsBattingExtra.Range("A16").Value = "'4290"
$wsBattingExtra.Range("B16").Value = 10
$wsBattingExtra.Range("F16").Value = "'NO"
$wsBattingExtra.Range("A17").Value = "'4378"
$wsBattingExtra.Range("B17").Value = 9
$wsBattingExtra.Range("C17").Value = "'1"
$wsBattingExtra.Range("D17").Value = "'0"
$wsBattingExtra.Range("E17").Value = "'6.50%"
$wsBattingExtra.Range("F17").Value = "'NO"
$wsBattingExtra.Range("A18").Value = "'4379"
$wsBattingExtra.Range("B18").Value = 9
$wsBattingExtra.Range("F18").Value = "'NO"
$wsBattingExtra.Range("A19").Value = "'4525"
$wsBattingExtra.Range("F19").Value = "'NO"
$wsBattingExtra.Range("A20").Value = "'4528"
$wsBattingExtra.Range("B20").Value = 8
$wsBattingExtra.Range("F20").Value = "'NO"

# --- "ODI Bowling Extra" header row ---
$wsBowlingExtra.Range("A1").Value = "'MATCH_CODE"
$wsBowlingExtra.Range("B1").Value = "'MAIDEN_OVERS"
$wsBowlingExtra.Range("C1").Value = "'PERCENT_WICKETS_OF_ALL"

# --- "ODI Bowling Extra" data rows ---
$wsBowlingExtra.Range("A2").Value = "'3649"
$wsBowlingExtra.Range("B2").Value = "'"
$wsBowlingExtra.Range("C2").Value = "'"
$wsBowlingExtra.Range("A3").Value = "'3650"
$wsBowlingExtra.Range("B3").Value = "'0"
$wsBowlingExtra.Range("C3").Value = "'10.00%"
$wsBowlingExtra.Range("A4").Value = "'3651"
$wsBowlingExtra.Range("B4").Value = "'"
$wsBowlingExtra.Range("C4").Value = "'"
$wsBowlingExtra.Range("A5").Value = "'3652"
$wsBowlingExtra.Range("B5").Value = "'0"
$wsBowlingExtra.Range("C5").Value = "'30.00%"
$wsBowlingExtra.Range("A6").Value = "'3707"
$wsBowlingExtra.Range("B6").Value = "'"
$wsBowlingExtra.Range("C6").Value = "'"
$wsBowlingExtra.Range("A7").Value = "'4134"
$wsBowlingExtra.Range("B7").Value = "'0"
$wsBowlingExtra.Range("C7").Value = "'10.00%"
$wsBowlingExtra.Range("A8").Value = "'4136"
$wsBowlingExtra.Range("B8").Value = "'0"
$wsBowlingExtra.Range("C8").Value = "'20.00%"
$wsBowlingExtra.Range("A9").Value = "'4140"
$wsBowlingExtra.Range("B9").Value = "'0"
$wsBowlingExtra.Range("C9").Value = "'"
$wsBowlingExtra.Range("A10").Value = "'4145"
$wsBowlingExtra.Range("B10").Value = "'0"
$wsBowlingExtra.Range("C10").Value = "'"
$wsBowlingExtra.Range("A11").Value = "'4147"
$wsBowlingExtra.Range("B11").Value = "'"
$wsBowlingExtra.Range("C11").Value = "'"
$wsBowlingExtra.Range("A12").Value = "'4154"
$wsBowlingExtra.Range("B12").Value = "'"
$wsBowlingExtra.Range("C12").Value = "'"
$wsBowlingExtra.Range("A13").Value = "'4160"
$wsBowlingExtra.Range("B13").Value = "'0"
$wsBowlingExtra.Range("C13").Value = "'"
$wsBowlingExtra.Range("A14").Value = "'4163"
$wsBowlingExtra.Range("B14").Value = "'1"
$wsBowlingExtra.Range("C14").Value = "'"
$wsBowlingExtra.Range("A15").Value = "'4164"
$wsBowlingExtra.Range("B15").Value = "'0"
$wsBowlingExtra.Range("C15").Value = "'10.00%"
$wsBowlingExtra.Range("A16").Value = "'4290"
$wsBowlingExtra.Range("B16").Value = "'0"
$wsBowlingExtra.Range("C16").Value = "'"
$wsBowlingExtra.Range("A17").Value = "'4378"
$wsBowlingExtra.Range("B17").Value = "'1"
$wsBowlingExtra.Range("C17").Value = "'10.00%"
$wsBowlingExtra.Range("A18").Value = "'4379"
$wsBowlingExtra.Range("B18").Value = "'0"
$wsBowlingExtra.Range("C18").Value = "'10.00%"
$wsBowlingExtra.Range("A19").Value = "'4525"
$wsBowlingExtra.Range("B19").Value = "'"
$wsBowlingExtra.Range("C19").Value = "'"
$wsBowlingExtra.Range("A20").Value = "'4528"
$wsBowlingExtra.Range("B20").Value = "'0"
$wsBowlingExtra.Range("C20").Value = "'"

# --- Copy the existing header formatting (bold, border, centered) onto the new headers ---
# (single-cell source so PasteSpecial tiles cleanly across the whole destination
#  row without spilling into an extra column)
$headerSource = $wsBatting.Range("A1")
[void]$headerSource.Copy()
[void]$wsBattingExtra.Range("A1:F1").PasteSpecial(-4122)
[void]$headerSource.Copy()
[void]$wsBowlingExtra.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore selection to A1 on each new sheet ---
[void]$wsBattingExtra.Range("A1").Select()
[void]$wsBowlingExtra.Range("A1").Select()

# --- Leave "Player Info" as the active sheet, matching the original workbook state ---
[void]$wb.Worksheets.Item("Player Info").Activate()

